$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows 2-47 (hourly refresh snapshot) ---
$ws.Range("D2").Value = "29.148.06"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.833.08"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'241.51"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'0.6636"
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.07436"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "'22.73"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").Value = "'0.07741"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "1.870.40"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "'4.992"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'0.6697"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "'83.01"
$ws.Range("E15").Value = "  -5.11%  "
$ws.Range("D16").Value = "'6.109"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "'0.000008370"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "29.179.84"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "'227.19"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").Value = "'12.49"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'7.169"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'159.83"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "'8.626"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D28").Value = "'1.512"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").Value = "'4.116"
$ws.Range("E29").Value = "  -3.23%  "
$ws.Range("D30").Value = "'4.045"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'0.05320"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "'1.871"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "'0.7531"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "'2.615"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "1.281.67"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").Value = "'0.01799"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "'0.9283"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "'0.08814"
$ws.Range("E41").Value = "  +11.28%  "
$ws.Range("D42").Value = "'5.960"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "'102.24"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("D45").Value = "1.983.08"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'0.5151"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "'1.769"
$ws.Range("E47").Value = "  -0.28%  "

# --- New coin BabyDogeCoin enters the ranking at #48, pushing Aave/Cronos/Aptos down one row; EnergySwap drops off the bottom ---
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'63.32"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05907"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.794"
$ws.Range("E51").Value = "  -1.39%  "
